# Refatorando o consolidador para modelo ETL
# Replace the sample absenteeism rows (2-11) with the new ETL-generated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  9776,  "Arthur das Neves",        "Jurídico",               "Viagem de negócios", 1, 45095, 4287.1),
    @(3,  18180, "Theo Gonçalves",          "Recursos Humanos",       "Consulta médica",    8, 45087, 8414.719999999999),
    @(4,  43768, "Maria Luiza Cavalcanti",  "Financeiro",              "Consulta médica",    7, 45081, 10456.49),
    @(5,  84237, "Nina Cavalcanti",         "Atendimento ao Cliente", "Consulta médica",    1, 45105, 3553.12),
    @(6,  13039, "Maria Clara Correia",     "Marketing",               "Viagem de negócios", 4, 45106, 9818.139999999999),
    @(7,  29096, "Caroline Cardoso",        "Atendimento ao Cliente", "Consulta médica",    2, 45106, 4872.14),
    @(8,  43354, "Maria Eduarda Oliveira",  "TI",                      "Outros",             2, 45105, 12213.18),
    @(9,  93774, "Benjamin Duarte",         "Recursos Humanos",       "Viagem de negócios", 1, 45100, 11423.31),
    @(10, 95453, "Luna Fogaça",             "Jurídico",               "Doença",             5, 45101, 9639.59),
    @(11, 47884, "Maria Barbosa",           "Operações",               "Viagem de negócios", 1, 45098, 12058.52)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
